# aggiornamento fino a 28/06 incluso
# Adds rows 270-301 (dates 2021-05-28 .. 2021-06-28) to Sheet1, extending
# the dimension from A1:D269 to A1:D301. New rows mirror the existing
# "no new cases" rows: date serial in col A (same style as prior date
# cells), and 0 in columns B, C, D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 269
$firstNewSerial = 44344
$lastNewRow = 301

# Use the last existing data row as the format template for the new rows.
$templateRow = $ws.Range("A$lastRow`:D$lastRow")

for ($r = $lastRow + 1; $r -le $lastNewRow; $r++) {
    $serial = $firstNewSerial + ($r - ($lastRow + 1))

    # Write values first.
    $ws.Range("A$r").Value = $serial
    $ws.Range("B$r").Value = 0
    $ws.Range("C$r").Value = 0
    $ws.Range("D$r").Value = 0

    # Copy formatting only (keeps the date number format / alignment /
    # border / bold font already used by column A, without disturbing the
    # values just written) from the template row.
    $templateRow.Copy()
    $ws.Range("A$r`:D$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
